# Generate Report for Handoff
# Updates the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" for the
# most recently handed-off file (fd002a3e-17cf-413a-a61f-4d1b2944eae1.md, row 7)
# across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G7").Value = "2016-08-25 16:44:13"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H7").Value = "2016-08-25 16:44:07"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H7").Value = "2016-08-25 16:44:13"
